# GitHub Actions crypto-price refresh
# Updates the Price (column D) and Volume(1h) (column E) cells on the
# coinranking.com-sourced rows of Sheet1 with freshly scraped values.
#
# Price values are written as text (NumberFormat "@") because the source
# data uses a dotted thousands/decimal notation (e.g. "58.784.49") that is
# not a valid Excel number, and we want every Price cell to keep the exact
# textual formatting that was scraped (trailing zeros, subscript digits, etc).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row number => @{ Price = <new price text, $null if unchanged>; Volume = <new Volume(1h) text> }
$updates = [ordered]@{
    2 = @{ Price = '58.784.49'; Volume = '  -1.94%  ' }
    3 = @{ Price = '2.303.85'; Volume = '  -4.27%  ' }
    4 = @{ Price = $null; Volume = '  -0.02%  ' }
    5 = @{ Price = '547.50'; Volume = '  -1.15%  ' }
    6 = @{ Price = '132.03'; Volume = '  -2.74%  ' }
    7 = @{ Price = $null; Volume = '  +0.01%  ' }
    8 = @{ Price = $null; Volume = '  -2.07%  ' }
    9 = @{ Price = '2.302.31'; Volume = '  -4.30%  ' }
    10 = @{ Price = $null; Volume = '  -2.87%  ' }
    11 = @{ Price = '5.52'; Volume = '  -1.54%  ' }
    12 = @{ Price = $null; Volume = '  +1.49%  ' }
    13 = @{ Price = $null; Volume = '  -4.62%  ' }
    14 = @{ Price = '23.88'; Volume = '  -3.06%  ' }
    15 = @{ Price = '2.713.72'; Volume = '  -4.35%  ' }
    16 = @{ Price = '58.772.67'; Volume = '  -1.74%  ' }
    18 = @{ Price = '2.320.21'; Volume = '  -3.61%  ' }
    19 = @{ Price = '10.69'; Volume = '  -4.25%  ' }
    20 = @{ Price = $null; Volume = '  -4.10%  ' }
    21 = @{ Price = '314.70'; Volume = '  -3.67%  ' }
    22 = @{ Price = '6.46'; Volume = '  -4.28%  ' }
    23 = @{ Price = $null; Volume = '  +0.05%  ' }
    24 = @{ Price = '63.49'; Volume = '  -1.73%  ' }
    25 = @{ Price = '0.169'; Volume = '  -5.98%  ' }
    26 = @{ Price = $null; Volume = '  +0.17%  ' }
    27 = @{ Price = '8.12'; Volume = '  -5.60%  ' }
    28 = @{ Price = '1.33'; Volume = '  -5.51%  ' }
    29 = @{ Price = $null; Volume = '  -1.72%  ' }
    30 = @{ Price = '168.75'; Volume = '  -0.86%  ' }
    31 = @{ Price = '0.0₃0725'; Volume = '  -5.10%  ' }
    32 = @{ Price = '1.10'; Volume = '  +1.36%  ' }
    33 = @{ Price = $null; Volume = '  -5.28%  ' }
    34 = @{ Price = '0.381'; Volume = '  -4.61%  ' }
    35 = @{ Price = $null; Volume = '  -0.02%  ' }
    36 = @{ Price = $null; Volume = '  -3.40%  ' }
    37 = @{ Price = $null; Volume = '  -0.03%  ' }
    38 = @{ Price = '1.26'; Volume = '  -4.19%  ' }
    39 = @{ Price = '3.98'; Volume = '  -4.95%  ' }
    40 = @{ Price = '38.09'; Volume = '  -1.18%  ' }
    41 = @{ Price = $null; Volume = '  -4.76%  ' }
    42 = @{ Price = '298.81'; Volume = '  -7.38%  ' }
    43 = @{ Price = '141.04'; Volume = '  -3.58%  ' }
    44 = @{ Price = '3.44'; Volume = '  -4.31%  ' }
    45 = @{ Price = $null; Volume = '  -1.07%  ' }
    46 = @{ Price = '0.0501'; Volume = '  -2.59%  ' }
    47 = @{ Price = '0.556'; Volume = '  -3.23%  ' }
    48 = @{ Price = '18.51'; Volume = '  -6.70%  ' }
    49 = @{ Price = $null; Volume = '  -2.41%  ' }
    50 = @{ Price = '16.64'; Volume = '  -3.49%  ' }
    51 = @{ Price = '11.04'; Volume = '  -0.05%  ' }
}

foreach ($row in $updates.Keys) {
    $change = $updates[$row]
    if ($null -ne $change.Price) {
        $priceCell = $ws.Range("D$row")
        $priceCell.NumberFormat = "@"
        $priceCell.Value = $change.Price
        $priceCell.Style = "Normal"
    }
    $ws.Range("E$row").Value = $change.Volume
}

